{"js": "// The target content replaces the whole body: the original single\n// underlined paragraph (\"\u00f3\u00f3\u00f3\u00f3\u00f3\") is gone, and is replaced by 16\n// paragraphs (three of which are fully empty) with no run/paragraph\n// formatting at all. We build the exact OOXML for the new body content\n// and inject it through Range.insertOoxml (Flat OPC), which lets us\n// control the precise markup (no stray empty runs on blank paragraphs,\n// no inherited rPr/pPr from the old underlined paragraph mark).\n\nconst paragraphs = [\n  \"Djajdlwjdliwjdl djaiwjdawidjaw\u00e7d\",\n  \"Jdlawdjawl pwajdpawj awdpoawjdpaow powadpo\",\n  \"Dawpodapowidapwodiawpodiaw[\",\n  \"\",\n  \"Dpoawidpoaiwd powaidpwaoidpaw\",\n  \"Ipodawid\u00b4powaid [awpoid powaidp[owa [poawidpo[aw powaid[poawi di[apw\",\n  \"Dapwodi[apow\",\n  \"\",\n  \"\",\n  \"Iwpodaiwd\",\n  \"Idpao\",\n  \"Adad \",\n  \"Adwdwadawi\",\n  \"A0wdiawpdia\",\n  \"Wdawidawdpiawpodi[\",\n  \"Owjdapowj awpd paw dpaw awpod apwod apowd iapow d\"\n];\n\nfunction escapeXml(text) {\n  return text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\")\n    .replace(/'/g, \"&apos;\");\n}\n\nfunction paragraphXml(text) {\n  if (text === \"\") {\n    return \"<w:p/>\";\n  }\n  const needsPreserve = /^\\s|\\s$/.test(text);\n  const preserve = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n  return `<w:p><w:r><w:t${preserve}>${escapeXml(text)}</w:t></w:r></w:p>`;\n}\n\nconst bodyParagraphsXml = paragraphs.map(paragraphXml).join(\"\");\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  bodyParagraphsXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst body = context.document.body;\n// Replace the whole body's content (this leaves the existing sectPr,\n// i.e. the page setup, untouched, same as the diff).\nbody.insertOoxml(flatOpcXml, \"Replace\");\n\nawait context.sync();\n", "ps1": "# The target content replaces the whole body: the original single\n# underlined paragraph (\"\u00f3\u00f3\u00f3\u00f3\u00f3\") is gone, and is replaced by 16\n# paragraphs (three of which are fully empty) with no run/paragraph\n# formatting at all. We build the exact OOXML for the new body content\n# and inject it through Range.InsertXML (Flat OPC), which lets us\n# control the precise markup (no stray empty runs on blank paragraphs,\n# no inherited rPr/pPr from the old underlined paragraph mark).\n\n$d = $word.ActiveDocument\n\n$paragraphs = @(\n  \"Djajdlwjdliwjdl djaiwjdawidjaw\u00e7d\",\n  \"Jdlawdjawl pwajdpawj awdpoawjdpaow powadpo\",\n  \"Dawpodapowidapwodiawpodiaw[\",\n  \"\",\n  \"Dpoawidpoaiwd powaidpwaoidpaw\",\n  \"Ipodawid\u00b4powaid [awpoid powaidp[owa [poawidpo[aw powaid[poawi di[apw\",\n  \"Dapwodi[apow\",\n  \"\",\n  \"\",\n  \"Iwpodaiwd\",\n  \"Idpao\",\n  \"Adad \",\n  \"Adwdwadawi\",\n  \"A0wdiawpdia\",\n  \"Wdawidawdpiawpodi[\",\n  \"Owjdapowj awpd paw dpaw awpod apwod apowd iapow d\"\n)\n\nfunction Escape-Xml([string]$text) {\n  return $text.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\").Replace('\"', \"&quot;\").Replace(\"'\", \"&apos;\")\n}\n\n$bodyXml = \"\"\nforeach ($p in $paragraphs) {\n  if ($p -eq \"\") {\n    $bodyXml += \"<w:p/>\"\n  } else {\n    $needsPreserve = ($p -match \"^\\s\") -or ($p -match \"\\s$\")\n    $preserveAttr = \"\"\n    if ($needsPreserve) { $preserveAttr = ' xml:space=\"preserve\"' }\n    $bodyXml += \"<w:p><w:r><w:t$preserveAttr>\" + (Escape-Xml $p) + \"</w:t></w:r></w:p>\"\n  }\n}\n\n$flatOpcXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  $bodyXml +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n# Replace the whole body's content (this leaves the existing sectPr,\n# i.e. the page setup, untouched, same as the diff).\n$d.Content.InsertXML($flatOpcXml)\n"}
